$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H15").Value = 244.5
$ws.Range("I15").Value = 244.5
$ws.Range("K15").Value = 733.5
$ws.Range("M15").Value = -564.5
$ws.Range("H18").Value = 1045.2106
$ws.Range("I18").Value = 755.7059
$ws.Range("J18").Value = 3506
$ws.Range("K18").Value = 755.7059
$ws.Range("L18").Value = 3506
$ws.Range("M18").Value = -471.7059
$ws.Range("N18").Value = -4074
$ws.Range("H40").Value = 4550
$ws.Range("J40").Value = 1200
$ws.Range("L40").Value = 1200
$ws.Range("N40").Value = -1550
$ws.Range("H44").Value = 17050
$ws.Range("J44").Value = 17050
$ws.Range("L44").Value = 17050
$ws.Range("N44").Value = -17974
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 15000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 15000
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -14514
$ws.Range("H138").Value = 2733.9546
$ws.Range("I138").Value = 1623.4615
$ws.Range("J138").Value = 3455.775
$ws.Range("K138").Value = 4870.3845
$ws.Range("L138").Value = 10367.325
$ws.Range("M138").Value = 269.6154999999999
$ws.Range("N138").Value = -20647.325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 7480
$ws.Range("H132").Value = 11704.225
$ws.Range("I132").Value = 13232.0205
$ws.Range("K132").Value = 39696.0615
$ws.Range("M132").Value = -37166.0615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 10833.333
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 9040.272000000001
$ws.Range("J19").Value = 48000
$ws.Range("L19").Value = 48000
$ws.Range("N19").Value = -48340
$ws.Range("H24").Value = 9040.272000000001
$ws.Range("J24").Value = 48000
$ws.Range("L24").Value = 48000
$ws.Range("N24").Value = -48340
$ws.Range("H132").Value = 55558364
$ws.Range("I132").Value = 71432056
$ws.Range("J132").Value = 38463624
$ws.Range("K132").Value = 214296168
$ws.Range("L132").Value = 115390872
$ws.Range("M132").Value = -214293638
$ws.Range("N132").Value = -115395932
$ws.Range("H134").Value = 1788
$ws.Range("I134").Value = 787.1429000000001
$ws.Range("K134").Value = 2361.4287
$ws.Range("M134").Value = 173.5712999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1094.3125
$ws.Range("I26").Value = 215.625
$ws.Range("J26").Value = 1973
$ws.Range("K26").Value = 646.875
$ws.Range("L26").Value = 5919
$ws.Range("M26").Value = -358.875
$ws.Range("N26").Value = -6495
$ws.Range("H70").Value = 68540.13
$ws.Range("I70").Value = 125962.75
$ws.Range("K70").Value = 377888.25
$ws.Range("M70").Value = -377573.25
$ws.Range("H73").Value = 68540.13
$ws.Range("I73").Value = 125962.75
$ws.Range("K73").Value = 377888.25
$ws.Range("M73").Value = -376796.25
$ws.Range("H80").Value = 1323.4445
$ws.Range("I80").Value = 833.3333
$ws.Range("J80").Value = 1421.4667
$ws.Range("K80").Value = 2499.9999
$ws.Range("L80").Value = 4264.4001
$ws.Range("M80").Value = -1563.9999
$ws.Range("N80").Value = -6136.4001
$ws.Range("H81").Value = 1926.6471
$ws.Range("I81").Value = 1016.1429
$ws.Range("J81").Value = 2564
$ws.Range("K81").Value = 3048.4287
$ws.Range("L81").Value = 7692
$ws.Range("M81").Value = -1925.4287
$ws.Range("N81").Value = -9938
$ws.Range("H83").Value = 1323.4445
$ws.Range("I83").Value = 833.3333
$ws.Range("J83").Value = 1421.4667
$ws.Range("K83").Value = 7499.9997
$ws.Range("L83").Value = 12793.2003
$ws.Range("M83").Value = -2819.9997
$ws.Range("N83").Value = -22153.2003
$ws.Range("H84").Value = 1926.6471
$ws.Range("I84").Value = 1016.1429
$ws.Range("J84").Value = 2564
$ws.Range("K84").Value = 9145.286100000001
$ws.Range("L84").Value = 23076
$ws.Range("M84").Value = -3529.286100000001
$ws.Range("N84").Value = -34308
$ws.Range("H137").Value = 14479936
$ws.Range("J137").Value = 27085752
$ws.Range("L137").Value = 81257256
$ws.Range("N137").Value = -81267456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 23331.334
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9369
$ws.Range("H70").Value = 93493.30499999999
$ws.Range("I70").Value = 172350.67
$ws.Range("J70").Value = 7467.091
$ws.Range("K70").Value = 172350.67
$ws.Range("L70").Value = 7467.091
$ws.Range("M70").Value = -172080.67
$ws.Range("N70").Value = -8007.091
$ws.Range("H73").Value = 93493.30499999999
$ws.Range("I73").Value = 172350.67
$ws.Range("J73").Value = 7467.091
$ws.Range("K73").Value = 172350.67
$ws.Range("L73").Value = 7467.091
$ws.Range("M73").Value = -171414.67
$ws.Range("N73").Value = -9339.091
$ws.Range("H132").Value = 2403.795
$ws.Range("I132").Value = 1677.2963
$ws.Range("J132").Value = 4038.4167
$ws.Range("K132").Value = 5031.8889
$ws.Range("L132").Value = 12115.2501
$ws.Range("M132").Value = -2501.8889
$ws.Range("N132").Value = -17175.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 558
$ws.Range("I30").Value = 558
$ws.Range("K30").Value = 558
$ws.Range("M30").Value = -450
$ws.Range("H57").Value = 7750
$ws.Range("J57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 35891
$ws.Range("I26").Value = 3836.5
$ws.Range("J26").Value = 100000
$ws.Range("K26").Value = 3836.5
$ws.Range("L26").Value = 100000
$ws.Range("M26").Value = -3543.5
$ws.Range("N26").Value = -100586
$ws.Range("H62").Value = 3848327.2
$ws.Range("I62").Value = 7694242
$ws.Range("J62").Value = 2412.2
$ws.Range("K62").Value = 7694242
$ws.Range("L62").Value = 2412.2
$ws.Range("M62").Value = -7693618
$ws.Range("N62").Value = -3660.2
$ws.Range("H65").Value = 3848327.2
$ws.Range("I65").Value = 7694242
$ws.Range("J65").Value = 2412.2
$ws.Range("K65").Value = 38471210
$ws.Range("L65").Value = 12061
$ws.Range("M65").Value = -38468090
$ws.Range("N65").Value = -18301
$ws.Range("H96").Value = 83334856
$ws.Range("I96").Value = 142858980
$ws.Range("J96").Value = 1078.8
$ws.Range("K96").Value = 142858980
$ws.Range("M96").Value = -142857607
$ws.Range("N96").Value = -3824.8
$ws.Range("H136").Value = 16185.865
$ws.Range("I136").Value = 29777.03
$ws.Range("J136").Value = 4633.375
$ws.Range("K136").Value = 89331.09
$ws.Range("L136").Value = 13900.125
$ws.Range("M136").Value = -86781.09
$ws.Range("N136").Value = -19000.125
